# Generate Report for Handoff
#
# Adds two newly-discovered localization files to the status report:
#   70cda61b-7cdf-4af1-b2fd-4d4f2d9d64d6.md  (row 4)
#   90d209b7-35b4-4236-8827-ded1df34bc6e.md  (row 5)
#
# Each file gets one row on the "Overview" sheet plus one row on each of the
# "zh-cn" and "de-de" per-language sheets (mirroring the existing rows for
# 7fd5598c-... / 114aac0a-... already present in the workbook).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Data describing the two new files.
# ---------------------------------------------------------------------------

$file1 = @{
    Uuid            = "70cda61b-7cdf-4af1-b2fd-4d4f2d9d64d6"
    OverviewDate    = "2016-32-18 22:32:12"
    ZhHash          = "b2a35c44c382a766a0410c7cdf8a37addd7b42de"
    ZhHandoffDate   = "2016-03-18 22:32:09"
    DeHash          = "b2a35c44c382a766a0410c7cdf8a37addd7b42de"
    DeHandoffDate   = "2016-03-18 22:32:12"
}

$file2 = @{
    Uuid            = "90d209b7-35b4-4236-8827-ded1df34bc6e"
    OverviewDate    = "2016-32-18 22:32:12"
    ZhHash          = "2d5be31e2a5a5b84a46e78a42e4c16f5fc9566a4"
    ZhHandoffDate   = "2016-03-18 22:32:09"
    DeHash          = "2d5be31e2a5a5b84a46e78a42e4c16f5fc9566a4"
    DeHandoffDate   = "2016-03-18 22:32:12"
}

$files = @($file1, $file2)

$readyForHandoff = "Ready for handoff"
$includeReason   = "Include"
$neverHandedBack = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRow = 4

foreach ($f in $files) {
    $mdName = "$($f.Uuid).md"

    $wsOverview.Hyperlinks.Add(
        $wsOverview.Range("A$overviewRow"),
        "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName",
        "",
        "",
        $mdName
    ) | Out-Null

    $wsOverview.Range("B$overviewRow").Value = $readyForHandoff
    $wsOverview.Range("C$overviewRow").Value = $readyForHandoff
    $wsOverview.Range("D$overviewRow").Value = $f.OverviewDate

    $overviewRow++
}

# ---------------------------------------------------------------------------
# Per-language sheets (zh-cn / de-de):
#   Source File Name | File Extension | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------

$languages = @(
    @{
        Sheet     = "zh-cn"
        HashProp  = "ZhHash"
        DateProp  = "ZhHandoffDate"
        RepoOrg   = "oltest.zh-cn"
    },
    @{
        Sheet     = "de-de"
        HashProp  = "DeHash"
        DateProp  = "DeHandoffDate"
        RepoOrg   = "oltest.de-de"
    }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)
    $row = 4

    foreach ($f in $files) {
        $mdName  = "$($f.Uuid).md"
        $hash    = $f[$lang.HashProp]
        $date    = $f[$lang.DateProp]
        $xlfName = "$($f.Uuid).$hash.$($lang.Sheet).xlf"

        # A: source file name (link to the .md source on github)
        $ws.Hyperlinks.Add(
            $ws.Range("A$row"),
            "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName",
            "",
            "",
            $mdName
        ) | Out-Null

        # B: file extension (displayed as ".md", links to the same source)
        $ws.Hyperlinks.Add(
            $ws.Range("B$row"),
            "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName",
            "",
            "",
            ".md"
        ) | Out-Null

        # C: status
        $ws.Range("C$row").Value = $readyForHandoff

        # D: latest handoff file (link to the generated .xlf handoff package)
        $ws.Hyperlinks.Add(
            $ws.Range("D$row"),
            "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/$($lang.RepoOrg)/ci/ht/$xlfName",
            "",
            "",
            $xlfName
        ) | Out-Null

        # E: latest handoff datetime
        $ws.Range("E$row").Value = $date
        $ws.Range("E$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"

        # F, G: latest target file / latest handback file -- left blank, the
        # file hasn't been handed back yet.

        # H: latest handback datetime -- never handed back
        $ws.Range("H$row").Value = $neverHandedBack

        # I: handoff reason
        $ws.Range("I$row").Value = $includeReason

        $row++
    }
}
